# Applies "another round of changes" commit to Config.xlsx
# - Settings sheet: replace the generic Framework/ProcessABCQueue example rows
#   with concrete Zoho queue/process/asset/url configuration rows, add a
#   hyperlink to the Zoho sign-in page, tweak the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- Row 2 : Orchestrator queue config -------------------------------------
# A2 / C2 keep their original text, only B2 (the example queue name) changes
# from "ProcessABCQueue" to "Zoho_Queue" and switches to the wrap-text style
# (matching the style already used by C4/C2-equivalent wrapped cells).
$ws.Range("B2").Value = "Zoho_Queue"
$ws.Range("B2").WrapText = $true

# --- Row 4 : Business process name config -----------------------------------
# B4 changes from the generic "Framework" example to "ZohoProcess".
$ws.Range("B4").Value = "ZohoProcess"

# --- Row 5 : Orchestrator asset config --------------------------------------
# A5 keeps "zoho_timesheet_orc". B5's text collapses from the 4x repeated
# "Zoho_Timesheet_Asset" lines down to a single value, vertically centered
# and still wrapped. C5 is a brand new description cell.
$ws.Range("B5").Value = "Zoho_Timesheet_Asset"
$ws.Range("B5").WrapText = $true
$ws.Range("B5").VerticalAlignment = -4108

$ws.Range("C5").Value = "Orchestrator asset Name. The value must match with the asset name defined on Orchestrator."

# --- Row 6 : brand new Zoho URL config row ----------------------------------
$ws.Range("A6").Value = "ZohoURL"
$ws.Range("C6").Value = "Zoho book sign in link "

$ws.Hyperlinks.Add($ws.Range("B6"), "https://accounts.zoho.com/signin", [Type]::Missing, [Type]::Missing, "https://accounts.zoho.com/signin")

# --- Active selection moves from B8 to C7 -----------------------------------
$ws.Activate()
$ws.Range("C7").Select()

$wb.Save()
